# Apply the reordered/updated "empadronador" monitoring data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2 through 15 (column A = name, column B = total_registros)
$names = @(
    "ZAMORA TAMAY NEYSER IVAN",
    "PÓSITO CHUGDEN NANIX",
    "TELLO FERNANDEZ MILENY",
    "VASQUEZ DIAZ LUZ ANGELICA",
    "ROJAS VASQUEZ FLOR NOELITA",
    "MEDINA VALLEJOS ERICK LEONARDO",
    "TIRADO PEREZ JEINER",
    "SOTO VILLENA NILSON",
    "SOTO VALLEJOS ELSITA",
    "RUIZ RUIZ LUZ MERI",
    "BENAVIDES MARRUFO ARACELI",
    "VASQUEZ LUNA YUDITH",
    "GALLARDO CORTEZ MELISSA DEL CARMEN",
    "BENAVIDES SALAZAR IDELSA"
)

$values = @(50, 49, 49, 49, 49, 48, 48, 47, 46, 44, 42, 42, 40, 37)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
